$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Increment the "days remaining" counter in column A for every existing
#    data row (rows 2..504) by 1, since a new day of data is appended below.
for ($r = 2; $r -le 504; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 1
}

# 2) Append the new day of data as row 505, copying formatting from the
#    last existing row (504) so styles (bold/border on A, date format on B)
#    carry over correctly.
$ws.Range("A504:L504").Copy()
$ws.Range("A505:L505").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(505, 1).Value2 = 0
$ws.Cells.Item(505, 2).Value2 = 43528
$ws.Cells.Item(505, 3).Value2 = 3845.09
$ws.Cells.Item(505, 4).Value2 = 3867.38
$ws.Cells.Item(505, 5).Value2 = 3733.75
$ws.Cells.Item(505, 6).Value2 = 3761.56
$ws.Cells.Item(505, 7).Value2 = 9029175788
$ws.Cells.Item(505, 8).Value2 = 66094551587
$ws.Cells.Item(505, 9).Value2 = -0.02250664740755148
$ws.Cells.Item(505, 10).Value2 = 0.02634605937350771
$ws.Cells.Item(505, 11).Value2 = 0.009049695223880027
$ws.Cells.Item(505, 12).Value2 = 5.414498483212308
